{"js": "// The document contains one table of two-digit x two-digit multiplication\n// problems (\"AA\u00d7BB=\") laid out five-per-row, with several blank spacer rows\n// in between. Every populated cell's text is replaced with a new problem,\n// left-to-right, top-to-bottom, in document order. Using Table.values lets\n// us rewrite just the <w:t> text while leaving all cell/paragraph/run\n// formatting (tcPr, pPr, rPr - fonts, size, etc.) untouched.\nconst oldToNew = [\n  [\"64\u00d779=\", \"43\u00d785=\"],\n  [\"88\u00d719=\", \"21\u00d794=\"],\n  [\"67\u00d722=\", \"56\u00d741=\"],\n  [\"26\u00d798=\", \"34\u00d758=\"],\n  [\"26\u00d798=\", \"37\u00d745=\"],\n  [\"24\u00d752=\", \"34\u00d753=\"],\n  [\"24\u00d761=\", \"32\u00d783=\"],\n  [\"39\u00d756=\", \"28\u00d793=\"],\n  [\"88\u00d798=\", \"63\u00d785=\"],\n  [\"81\u00d735=\", \"17\u00d715=\"],\n  [\"40\u00d737=\", \"73\u00d774=\"],\n  [\"79\u00d799=\", \"59\u00d769=\"],\n  [\"69\u00d772=\", \"81\u00d793=\"],\n  [\"50\u00d788=\", \"37\u00d711=\"],\n  [\"32\u00d744=\", \"55\u00d738=\"],\n  [\"38\u00d769=\", \"90\u00d742=\"],\n  [\"26\u00d785=\", \"37\u00d744=\"],\n  [\"23\u00d738=\", \"86\u00d785=\"],\n  [\"62\u00d783=\", \"50\u00d775=\"],\n  [\"77\u00d723=\", \"62\u00d783=\"],\n  [\"78\u00d713=\", \"25\u00d749=\"],\n  [\"42\u00d793=\", \"83\u00d715=\"],\n  [\"18\u00d746=\", \"14\u00d727=\"],\n  [\"77\u00d741=\", \"50\u00d773=\"],\n  [\"72\u00d773=\", \"99\u00d769=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet pairIdx = 0;\nfor (const table of tables.items) {\n  table.load(\"values\");\n  await context.sync();\n\n  const values = table.values;\n  for (let r = 0; r < values.length; r++) {\n    const row = values[r];\n    for (let c = 0; c < row.length; c++) {\n      const cellText = row[c];\n      if (cellText === \"\" || cellText === null || cellText === undefined) {\n        continue;\n      }\n      if (pairIdx >= oldToNew.length) {\n        continue;\n      }\n      const [expectedOld, replacement] = oldToNew[pairIdx];\n      if (cellText !== expectedOld) {\n        throw new Error(\n          \"Unexpected cell text '\" + cellText + \"' at row \" + r +\n          \", col \" + c + \"; expected '\" + expectedOld + \"'\"\n        );\n      }\n      row[c] = replacement;\n      pairIdx++;\n    }\n  }\n  table.values = values;\n  await context.sync();\n}\n", "ps1": "# The document contains one table of two-digit x two-digit multiplication\n# problems (\"AA\u00d7BB=\") laid out five-per-row, with several blank spacer rows\n# in between. Every populated cell's text is replaced with a new problem,\n# left-to-right, top-to-bottom, in document order. Writing directly to\n# Cell.Range.Text only rewrites the text run content, leaving cell/paragraph/\n# run formatting (tcPr, pPr, rPr - fonts, size, etc.) untouched.\n\n$oldValues = @(\n  \"64\u00d779=\", \"88\u00d719=\", \"67\u00d722=\", \"26\u00d798=\", \"26\u00d798=\",\n  \"24\u00d752=\", \"24\u00d761=\", \"39\u00d756=\", \"88\u00d798=\", \"81\u00d735=\",\n  \"40\u00d737=\", \"79\u00d799=\", \"69\u00d772=\", \"50\u00d788=\", \"32\u00d744=\",\n  \"38\u00d769=\", \"26\u00d785=\", \"23\u00d738=\", \"62\u00d783=\", \"77\u00d723=\",\n  \"78\u00d713=\", \"42\u00d793=\", \"18\u00d746=\", \"77\u00d741=\", \"72\u00d773=\"\n)\n$newValues = @(\n  \"43\u00d785=\", \"21\u00d794=\", \"56\u00d741=\", \"34\u00d758=\", \"37\u00d745=\",\n  \"34\u00d753=\", \"32\u00d783=\", \"28\u00d793=\", \"63\u00d785=\", \"17\u00d715=\",\n  \"73\u00d774=\", \"59\u00d769=\", \"81\u00d793=\", \"37\u00d711=\", \"55\u00d738=\",\n  \"90\u00d742=\", \"37\u00d744=\", \"86\u00d785=\", \"50\u00d775=\", \"62\u00d783=\",\n  \"25\u00d749=\", \"83\u00d715=\", \"14\u00d727=\", \"50\u00d773=\", \"99\u00d769=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$pairIdx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cellRange = $cell.Range\n    # Cell.Range.Text includes a trailing cell-mark character; trim it off\n    # so we compare/match on the visible text only.\n    $txt = $cellRange.Text\n    $txt = $txt.TrimEnd([char]7, [char]13)\n    if ($txt.Length -eq 0) {\n      continue\n    }\n    if ($pairIdx -ge $oldValues.Length) {\n      continue\n    }\n    $expectedOld = $oldValues[$pairIdx]\n    if ($txt -ne $expectedOld) {\n      throw \"Unexpected cell text '$txt' at row $r, col $c; expected '$expectedOld'\"\n    }\n    $cellRange.Text = $newValues[$pairIdx]\n    $pairIdx++\n  }\n}\n"}
